# Add "NA" values in column E (duplicate_image_filename) for rows 2-21,
# matching the rest of the stimuli data block (rows 2-21 under the header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
